$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "@"
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "69.794.34"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "3.498.70"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "580.81"
$ws.Range("E5").Value = "  -4.45%  "
Set-TextValue $ws.Range("D6") "193.66"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("D8").Value = "3.484.00"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("E9").Value = "  +0.04%  "
Set-TextValue $ws.Range("D10") "0.205"
$ws.Range("E10").Value = "  -7.08%  "
Set-TextValue $ws.Range("D11") "0.620"
$ws.Range("E11").Value = "  -4.29%  "
Set-TextValue $ws.Range("D12") "51.64"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("E13").Value = "  -6.21%  "
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "4.056.38"
$ws.Range("E15").Value = "  -3.52%  "
Set-TextValue $ws.Range("D16") "649.64"
$ws.Range("E16").Value = "  -5.24%  "
$ws.Range("D17").Value = "69.703.47"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "3.494.91"
$ws.Range("E18").Value = "  -5.45%  "
Set-TextValue $ws.Range("D19") "12.38"
$ws.Range("E19").Value = "  -4.37%  "
$ws.Range("E20").Value = "  -1.74%  "
Set-TextValue $ws.Range("D21") "18.28"
$ws.Range("E21").Value = "  -3.95%  "
Set-TextValue $ws.Range("D22") "0.950"
$ws.Range("E22").Value = "  -4.90%  "
Set-TextValue $ws.Range("D23") "18.10"
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("E24").Value = "  -2.68%  "
Set-TextValue $ws.Range("D25") "98.88"
$ws.Range("E25").Value = "  -6.13%  "
$ws.Range("E26").Value = "  -7.49%  "
$ws.Range("E27").Value = "  -3.67%  "
$ws.Range("E28").Value = "  -4.03%  "
Set-TextValue $ws.Range("D29") "9.34"
$ws.Range("E29").Value = "  -6.34%  "
Set-TextValue $ws.Range("D30") "32.70"
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("E31").Value = "  -6.44%  "
$ws.Range("E32").Value = "  -6.19%  "
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("E34").Value = "  -4.71%  "
Set-TextValue $ws.Range("D35") "61.19"
$ws.Range("E35").Value = "  -3.29%  "
Set-TextValue $ws.Range("D36") "536.87"
$ws.Range("E36").Value = "  +7.30%  "
$ws.Range("D37").Value = "3.704.16"
$ws.Range("E37").Value = "  -6.43%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "0.0₃0790"
$ws.Range("E39").Value = "  -8.99%  "
Set-TextValue $ws.Range("D40") "3.55"
$ws.Range("E40").Value = "  -0.04%  "
Set-TextValue $ws.Range("D41") "2.92"
$ws.Range("E41").Value = "  -4.04%  "
Set-TextValue $ws.Range("D42") "0.374"
$ws.Range("E42").Value = "  -3.67%  "
Set-TextValue $ws.Range("D43") "3.55"
$ws.Range("E43").Value = "  +41.53%  "
$ws.Range("E44").Value = "  -2.41%  "
Set-TextValue $ws.Range("D45") "34.37"
$ws.Range("E45").Value = "  -6.49%  "
$ws.Range("E46").Value = "  -2.83%  "
Set-TextValue $ws.Range("D47") "3.39"
$ws.Range("E47").Value = "  -2.91%  "
Set-TextValue $ws.Range("D48") "2.83"
$ws.Range("E48").Value = "  -7.31%  "
$ws.Range("E49").Value = "  -4.12%  "
Set-TextValue $ws.Range("D50") "0.999"
$ws.Range("E50").Value = "  -0.30%  "
Set-TextValue $ws.Range("D51") "8.21"
$ws.Range("E51").Value = "  -5.23%  "
